{"js": "// Add a second exercise link (\"ex2\") after the existing \"ex1\" hyperlink in\n// the \"LATIHAN DOM ... latihanWarna  latihanGameSuwitJawa ... ex1\" paragraph.\n//\n// Target XML shape (see diff):\n//   ...</w:hyperlink>            <- closes the \"ex1\" hyperlink\n//   <w:r><w:t xml:space=\"preserve\">  </w:t></w:r>\n//   <w:hyperlink ...><w:r><w:rPr><w:rStyle w:val=\"Hyperlink\"/></w:rPr><w:t>ex2</w:t></w:r></w:hyperlink>\n//   </w:p>\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// Find the paragraph that contains the \"ex1\" exercise link we need to extend.\n// It is the paragraph that starts with \"LATIHAN DOM\" and also contains\n// \"latihanGameSuwitJawa\" (there are two \"LATIHAN DOM\" paragraphs in the doc;\n// this is the one with the \"ex1\" link, not \"buatImageGallery\").\nlet targetParagraph = null;\nfor (const para of paragraphs.items) {\n  if (para.text.indexOf(\"latihanGameSuwitJawa\") !== -1) {\n    targetParagraph = para;\n    break;\n  }\n}\n\nif (!targetParagraph) {\n  throw new Error(\"Could not locate the 'LATIHAN DOM' / latihanGameSuwitJawa paragraph\");\n}\n\n// Insert \"  \" (two spaces) at the very end of the paragraph (before the\n// paragraph mark), matching the two-space separator used between the other\n// hyperlinks in this paragraph.\nconst spacerRange = targetParagraph.getRange(\"End\");\nspacerRange.insertText(\"  \", \"End\");\nawait context.sync();\n\n// Insert the new hyperlink text right after the spacer, then turn the\n// inserted text into a real hyperlink (this produces a <w:hyperlink> with\n// a run styled with the built-in \"Hyperlink\" character style, exactly like\n// the existing links in this document).\nconst linkTextRange = targetParagraph.getRange(\"End\");\nlinkTextRange.insertText(\"ex2\", \"End\");\nawait context.sync();\n\n// Re-locate the freshly inserted \"ex2\" text (the last occurrence in the\n// paragraph) and assign its hyperlink target.\nconst matches = targetParagraph.search(\"ex2\", { matchCase: true });\nmatches.load(\"items\");\nawait context.sync();\n\nconst newLinkRange = matches.items[matches.items.length - 1];\nnewLinkRange.hyperlink = \"ex/latihangamesuwitjawa.htm\";\nawait context.sync();\n", "ps1": "# Add a second exercise link (\"ex2\") after the existing \"ex1\" hyperlink in\n# the \"LATIHAN DOM ... latihanWarna  latihanGameSuwitJawa ... ex1\" paragraph.\n#\n# Target XML shape (see diff):\n#   ...</w:hyperlink>            <- closes the \"ex1\" hyperlink\n#   <w:r><w:t xml:space=\"preserve\">  </w:t></w:r>\n#   <w:hyperlink ...><w:r><w:rPr><w:rStyle w:val=\"Hyperlink\"/></w:rPr><w:t>ex2</w:t></w:r></w:hyperlink>\n#   </w:p>\n\n$d = $word.ActiveDocument\n\n# Find the paragraph that contains the \"ex1\" exercise link we need to extend.\n# It is the paragraph that contains \"latihanGameSuwitJawa\" (there are two\n# \"LATIHAN DOM\" paragraphs in the doc; this is the one with the \"ex1\" link,\n# not the later one with \"buatImageGallery\").\n$targetPara = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*latihanGameSuwitJawa*\") {\n        $targetPara = $p\n        break\n    }\n}\n\n# The paragraph's Range.End sits just past the paragraph mark, so back up\n# one character to get the insertion point right before it.\n$endPos = $targetPara.Range.End - 1\n$insertionRange = $d.Range($endPos, $endPos)\n\n# Insert the two-space separator, then the \"ex2\" display text, as plain\n# text first (matches the separator style used elsewhere in this paragraph).\n$insertionRange.InsertAfter(\"  ex2\")\n\n# Re-locate the \"ex2\" text we just inserted (it is now the text immediately\n# before the paragraph mark) and convert it into a real hyperlink.\n$linkEnd = $targetPara.Range.End - 1\n$linkStart = $linkEnd - 3\n$linkRange = $d.Range($linkStart, $linkEnd)\n\n$d.Hyperlinks.Add($linkRange, \"ex/latihangamesuwitjawa.htm\", \"\", \"\", \"ex2\") | Out-Null\n"}
